$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy header formatting (bold/border/center) from D1:E1 into F1:I1 ---
$ws.Range("D1:E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("H1:I1").PasteSpecial(-4122)

# --- New header labels ---
$ws.Range("F1").Value = "09-04-2025 Status"
$ws.Range("G1").Value = "09-04-2025 Time"
$ws.Range("H1").Value = "10-04-2025 Status"
$ws.Range("I1").Value = "10-04-2025 Time"

# --- Fill the new attendance columns for every student row (2-22) ---
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 6).Value = "A"
    $ws.Cells.Item($r, 7).Value = "00:00:00"
    $ws.Cells.Item($r, 8).Value = "A"
    $ws.Cells.Item($r, 9).Value = "00:00:00"
}

# --- Row 8 (ST3007 / ESA YOUSAF) attendance update for the 07-04-2025 columns ---
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "P"
$ws.Range("E8").Value = "09:46:06 AM"
